# Scheduled runner update: refresh market-board price/profit figures
# (currentAveragePrice*, Leve Price*, Leve Profit*) across several
# crafting-job leve sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 168
$ws.Range("I12").Value = 98.75
$ws.Range("K12").Value = 98.75
$ws.Range("M12").Value = 71.25

$ws.Range("H15").Value = 754.1754
$ws.Range("I15").Value = 754.1754
$ws.Range("K15").Value = 2262.5262
$ws.Range("M15").Value = -2093.5262

$ws.Range("H53").Value = 539.73334
$ws.Range("I53").Value = 161.2
$ws.Range("K53").Value = 161.2
$ws.Range("M53").Value = 475.8

$ws.Range("H106").Value = 2000
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H107").Value = 1225.6111
$ws.Range("I107").Value = 1071.1333
$ws.Range("K107").Value = 1071.1333
$ws.Range("M107").Value = 848.8667

$ws.Range("H132").Value = 1430.6111
$ws.Range("I132").Value = 1440.9062
$ws.Range("J132").Value = 1348.25
$ws.Range("K132").Value = 4322.7186
$ws.Range("L132").Value = 4044.75
$ws.Range("M132").Value = -1792.7186
$ws.Range("N132").Value = -9104.75

$ws.Range("H138").Value = 2313.205
$ws.Range("I138").Value = 1655
$ws.Range("J138").Value = 2938.5
$ws.Range("K138").Value = 4965
$ws.Range("L138").Value = 8815.5
$ws.Range("M138").Value = 175
$ws.Range("N138").Value = -19095.5

$ws.Range("H141").Value = 2488.36
$ws.Range("I141").Value = 2282.318
$ws.Range("K141").Value = 6846.954000000001
$ws.Range("M141").Value = -1666.954000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 499
$ws.Range("J12").Value = 499
$ws.Range("L12").Value = 499
$ws.Range("N12").Value = -845

$ws.Range("H32").Value = 5440.0166
$ws.Range("I32").Value = 2439.204
$ws.Range("K32").Value = 2439.204
$ws.Range("M32").Value = -2152.204

$ws.Range("H45").Value = 8408677
$ws.Range("I45").Value = 2265.111
$ws.Range("K45").Value = 2265.111
$ws.Range("M45").Value = -1888.111

$ws.Range("H61").Value = 58751.61
$ws.Range("I61").Value = 3381.8333
$ws.Range("K61").Value = 3381.8333
$ws.Range("M61").Value = -3169.8333

$ws.Range("H74").Value = 4902.3335
$ws.Range("I74").Value = 2480.4546
$ws.Range("J74").Value = 7566.4
$ws.Range("K74").Value = 2480.4546
$ws.Range("L74").Value = 7566.4
$ws.Range("M74").Value = -1606.4546
$ws.Range("N74").Value = -9314.4

$ws.Range("H77").Value = 4902.3335
$ws.Range("I77").Value = 2480.4546
$ws.Range("J77").Value = 7566.4
$ws.Range("K77").Value = 12402.273
$ws.Range("L77").Value = 37832
$ws.Range("M77").Value = -8034.273000000001
$ws.Range("N77").Value = -46568

$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws.Range("H132").Value = 3306.875
$ws.Range("I132").Value = 2537.5757
$ws.Range("J132").Value = 6933.5713
$ws.Range("K132").Value = 7612.7271
$ws.Range("L132").Value = 20800.7139
$ws.Range("M132").Value = -5082.7271
$ws.Range("N132").Value = -25860.7139

$ws.Range("H136").Value = 58751.61
$ws.Range("I136").Value = 3381.8333
$ws.Range("K136").Value = 10145.4999
$ws.Range("M136").Value = -7595.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 59338.668
$ws.Range("J13").Value = 59338.668
$ws.Range("L13").Value = 59338.668
$ws.Range("N13").Value = -59674.668

$ws.Range("H31").Value = 3250
$ws.Range("J31").Value = 3250
$ws.Range("L31").Value = 3250
$ws.Range("N31").Value = -3754

$ws.Range("H105").Value = 256612
$ws.Range("I105").Value = 510000
$ws.Range("K105").Value = 510000
$ws.Range("M105").Value = -508253

$ws.Range("H134").Value = 5706.0967
$ws.Range("I134").Value = 3444.6667
$ws.Range("J134").Value = 10455.1
$ws.Range("K134").Value = 10334.0001
$ws.Range("L134").Value = 31365.3
$ws.Range("M134").Value = -7799.000100000001
$ws.Range("N134").Value = -36435.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2249.5151
$ws.Range("I31").Value = 1427.4706
$ws.Range("J31").Value = 3122.9375
$ws.Range("K31").Value = 1427.4706
$ws.Range("L31").Value = 3122.9375
$ws.Range("M31").Value = -1132.4706
$ws.Range("N31").Value = -3712.9375

$ws.Range("H34").Value = 2249.5151
$ws.Range("I34").Value = 1427.4706
$ws.Range("J34").Value = 3122.9375
$ws.Range("K34").Value = 1427.4706
$ws.Range("L34").Value = 3122.9375
$ws.Range("M34").Value = -1225.4706
$ws.Range("N34").Value = -3526.9375

$ws.Range("H58").Value = 3216.8
$ws.Range("J58").Value = 3486.9092
$ws.Range("L58").Value = 3486.9092
$ws.Range("N58").Value = -3892.9092

$ws.Range("H107").Value = 1002.8261
$ws.Range("I107").Value = 729.2941
$ws.Range("J107").Value = 1777.8334
$ws.Range("K107").Value = 729.2941
$ws.Range("L107").Value = 1777.8334
$ws.Range("M107").Value = 1190.7059
$ws.Range("N107").Value = -5617.8334

$ws.Range("H122").Value = 4941.8184
$ws.Range("J122").Value = 5469.8
$ws.Range("L122").Value = 16409.4
$ws.Range("N122").Value = -21309.4

$ws.Range("H132").Value = 1395785
$ws.Range("I132").Value = 1381831.1
$ws.Range("K132").Value = 4145493.3
$ws.Range("M132").Value = -4142963.3

$ws.Range("H134").Value = 11910132
$ws.Range("I134").Value = 23816930
$ws.Range("J134").Value = 3333.3333
$ws.Range("K134").Value = 71450790
$ws.Range("L134").Value = 9999.999899999999
$ws.Range("M134").Value = -71448255
$ws.Range("N134").Value = -15069.9999

$ws.Range("H136").Value = 3216.8
$ws.Range("J136").Value = 3486.9092
$ws.Range("L136").Value = 10460.7276
$ws.Range("N136").Value = -15560.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1639.4445
$ws.Range("I131").Value = 1197
$ws.Range("K131").Value = 3591
$ws.Range("M131").Value = 1449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9898.875
$ws.Range("I132").Value = 5233.5
$ws.Range("J132").Value = 14564.25
$ws.Range("K132").Value = 15700.5
$ws.Range("L132").Value = 43692.75
$ws.Range("M132").Value = -13170.5
$ws.Range("N132").Value = -48752.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 1000
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("N26").Value = -1590

$ws.Range("H93").Value = 2253.5715
$ws.Range("I93").Value = 2253.5715
$ws.Range("K93").Value = 2253.5715
$ws.Range("M93").Value = -1005.5715

$ws.Range("H132").Value = 4337.2
$ws.Range("I132").Value = 3147.8572
$ws.Range("K132").Value = 9443.571599999999
$ws.Range("M132").Value = -6913.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 11139.8
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 11139.8
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 11139.8
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -11835.8

$ws.Range("H40").Value = 9000
$ws.Range("J40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("N40").Value = -9298

$ws.Range("H96").Value = 4388596.5
$ws.Range("J96").Value = 7520824
$ws.Range("L96").Value = 7520824
$ws.Range("N96").Value = -7523570

$ws.Range("H132").Value = 2531.7073
$ws.Range("I132").Value = 2251.4856
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 6754.4568
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -4224.4568
$ws.Range("N132").Value = -17559.0005

$ws.Range("H136").Value = 2730
$ws.Range("I136").Value = 2594.4285
$ws.Range("J136").Value = 3204.5
$ws.Range("K136").Value = 7783.2855
$ws.Range("L136").Value = 9613.5
$ws.Range("M136").Value = -5233.2855
$ws.Range("N136").Value = -14713.5
